$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 10 (A10 = 8) : HKAP2 Seq1.1 Hits 4 / Hong Kong Action Kit - Hit Kit #2 /
#                    NoizBoy / Objective: Assassination
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "HKAP2 Seq1.1 Hits 4"
$ws.Range("C10").NumberFormat = "h:mm"
$ws.Range("C10").Value = 0.00069444444444444447
$ws.Range("D10").Value = "Hong Kong Action Kit - Hit Kit #2"
$ws.Range("E10").Value = "NoizBoy"
$ws.Range("F10").Value = "Objective: Assassination"
$ws.Range("G10").Value = "Sonniss.com - GDC 2015 - Game Audio Bundle"

# ---------------------------------------------------------------------------
# Row 8 (A8 = 6) : Text Type-Hi Tech-003 / Modular UI /
#                  SoundMorph - Richard Devine / Objective: Bomb Defusing
# ---------------------------------------------------------------------------
$ws.Range("E8").Value = "SoundMorph - Richard Devine"
$ws.Range("D8").Value = "Modular UI"
$ws.Range("B8").Value = "Text Type-Hi Tech-003"
$ws.Range("C8").NumberFormat = "h:mm"
$ws.Range("C8").Value = 0.0090277777777777787
$ws.Range("F8").Value = "Objective: Bomb Defusing"
$ws.Range("G8").Value = "Sonniss.com - GDC 2015 - Game Audio Bundle"
# B8 loses its right border in the authored edit (new border/style pair)
$ws.Range("B8").Borders.Item(10).LineStyle = -4142

# ---------------------------------------------------------------------------
# Row 9 (A9 = 7) : Data Transmissions-001 / Modular UI /
#                  SoundMorph - Richard Devine / Objective: Bomb Planting
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = "Data Transmissions-001"
$ws.Range("C9").NumberFormat = "h:mm"
$ws.Range("C9").Value = 0.0097222222222222224
$ws.Range("D9").Value = "Modular UI"
$ws.Range("E9").Value = "SoundMorph - Richard Devine"
$ws.Range("F9").Value = "Objective: Bomb Planting"
$ws.Range("G9").Value = "Sonniss.com - GDC 2015 - Game Audio Bundle"

# ---------------------------------------------------------------------------
# Row 11 (A11 = 9) : EFX EXT .50 Cal Pistol Shots 02 A / Guns /
#                    Coll Anderson / Death
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "EFX EXT .50 Cal Pistol Shots 02 A"
$ws.Range("C11").NumberFormat = "h:mm"
$ws.Range("C11").Value = 0.0055555555555555558
$ws.Range("D11").Value = "Guns"
$ws.Range("E11").Value = "Coll Anderson"
$ws.Range("F11").Value = "Death"
$ws.Range("G11").Value = "Sonniss.com - GDC 2015 - Game Audio Bundle"

# ---------------------------------------------------------------------------
# Update selection to match the authored edit
# ---------------------------------------------------------------------------
$ws.Range("F13").Select()
